# Horarios actualizados Linea 141 - 888
# The scraper re-ran (Hora_Scrap 07:38:30) and the per-sheet schedules were
# re-sorted by Hora_Llegada; this updates every cell that moved/changed and
# appends the newly scraped rows at the bottom of each sheet's data block.
$wb = $excel.ActiveWorkbook

# ----- LP1912 (31 changed rows) -----
$ws = $wb.Worksheets.Item('LP1912')
$ws.Cells.Item(2, 1).Value = 'Última actualización: 07:38:30'
$ws.Cells.Item(3, 1).Value = 'Total filas: 92'
$ws.Cells.Item(49, 1).Value = '06:46:37'
$ws.Cells.Item(49, 3).Value = '14_ABASTO'
$ws.Cells.Item(49, 4).Value = 8
$ws.Cells.Item(50, 1).Value = '06:53:56'
$ws.Cells.Item(50, 3).Value = '17_ROMERO'
$ws.Cells.Item(50, 4).Value = 1
$ws.Cells.Item(65, 1).Value = '06:46:37'
$ws.Cells.Item(65, 3).Value = '17X38_ROMERO'
$ws.Cells.Item(65, 4).Value = 50
$ws.Cells.Item(66, 1).Value = '07:12:47'
$ws.Cells.Item(66, 3).Value = '27_EL RETIRO'
$ws.Cells.Item(66, 4).Value = 24
$ws.Cells.Item(71, 1).Value = '07:38:30'
$ws.Cells.Item(71, 4).Value = 11
$ws.Cells.Item(72, 1).Value = '07:38:30'
$ws.Cells.Item(72, 4).Value = 20
$ws.Cells.Item(73, 1).Value = '07:38:30'
$ws.Cells.Item(73, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(73, 4).Value = 21
$ws.Cells.Item(74, 1).Value = '06:53:56'
$ws.Cells.Item(74, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(74, 4).Value = 66
$ws.Cells.Item(75, 1).Value = '06:46:37'
$ws.Cells.Item(75, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(75, 4).Value = 74
$ws.Cells.Item(76, 1).Value = '06:18:01'
$ws.Cells.Item(76, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(76, 4).Value = 102
$ws.Cells.Item(78, 1).Value = '07:38:30'
$ws.Cells.Item(78, 4).Value = 25
$ws.Cells.Item(80, 1).Value = '07:38:30'
$ws.Cells.Item(80, 4).Value = 36
$ws.Cells.Item(81, 1).Value = '07:38:30'
$ws.Cells.Item(81, 4).Value = 41
$ws.Cells.Item(82, 1).Value = '07:38:30'
$ws.Cells.Item(82, 2).Value = '08:21'
$ws.Cells.Item(82, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(82, 4).Value = 43
$ws.Cells.Item(83, 1).Value = '07:38:30'
$ws.Cells.Item(83, 2).Value = '08:29'
$ws.Cells.Item(83, 3).Value = '14_ABASTO'
$ws.Cells.Item(83, 4).Value = 51
$ws.Cells.Item(84, 1).Value = '07:12:47'
$ws.Cells.Item(84, 2).Value = '08:33'
$ws.Cells.Item(84, 4).Value = 81
$ws.Cells.Item(85, 1).Value = '07:38:30'
$ws.Cells.Item(85, 2).Value = '08:34'
$ws.Cells.Item(85, 3).Value = '215C_EL PATO'
$ws.Cells.Item(85, 4).Value = 56
$ws.Cells.Item(86, 1).Value = '07:38:30'
$ws.Cells.Item(86, 2).Value = '08:41'
$ws.Cells.Item(86, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(86, 4).Value = 63
$ws.Cells.Item(87, 2).Value = '08:47'
$ws.Cells.Item(87, 3).Value = '215A_EL PATO'
$ws.Cells.Item(87, 4).Value = 95
# new row 88 (freshly scraped)
$ws.Cells.Item(88, 1).Value = '07:38:30'
$ws.Cells.Item(88, 2).Value = '08:48'
$ws.Cells.Item(88, 3).Value = '215A_EL PATO'
$ws.Cells.Item(88, 4).Value = 70
$ws.Cells.Item(88, 5).Value = 'LP1912'
# new row 89 (freshly scraped)
$ws.Cells.Item(89, 1).Value = '07:38:30'
$ws.Cells.Item(89, 2).Value = '08:51'
$ws.Cells.Item(89, 3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(89, 4).Value = 73
$ws.Cells.Item(89, 5).Value = 'LP1912'
# new row 90 (freshly scraped)
$ws.Cells.Item(90, 1).Value = '07:38:30'
$ws.Cells.Item(90, 2).Value = '08:59'
$ws.Cells.Item(90, 3).Value = '215B_EL PATO'
$ws.Cells.Item(90, 4).Value = 81
$ws.Cells.Item(90, 5).Value = 'LP1912'
# new row 91 (freshly scraped)
$ws.Cells.Item(91, 1).Value = '07:38:30'
$ws.Cells.Item(91, 2).Value = '09:02'
$ws.Cells.Item(91, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(91, 4).Value = 84
$ws.Cells.Item(91, 5).Value = 'LP1912'
# new row 92 (freshly scraped)
$ws.Cells.Item(92, 1).Value = '07:38:30'
$ws.Cells.Item(92, 2).Value = '09:03'
$ws.Cells.Item(92, 3).Value = '17X38_ROMERO'
$ws.Cells.Item(92, 4).Value = 85
$ws.Cells.Item(92, 5).Value = 'LP1912'
# new row 93 (freshly scraped)
$ws.Cells.Item(93, 1).Value = '07:38:30'
$ws.Cells.Item(93, 2).Value = '09:14'
$ws.Cells.Item(93, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(93, 4).Value = 96
$ws.Cells.Item(93, 5).Value = 'LP1912'
# new row 94 (freshly scraped)
$ws.Cells.Item(94, 1).Value = '07:38:30'
$ws.Cells.Item(94, 2).Value = '09:15'
$ws.Cells.Item(94, 3).Value = '27_EL RETIRO'
$ws.Cells.Item(94, 4).Value = 97
$ws.Cells.Item(94, 5).Value = 'LP1912'
# new row 95 (freshly scraped)
$ws.Cells.Item(95, 1).Value = '07:38:30'
$ws.Cells.Item(95, 2).Value = '09:19'
$ws.Cells.Item(95, 3).Value = '215_EL PELIGRO'
$ws.Cells.Item(95, 4).Value = 101
$ws.Cells.Item(95, 5).Value = 'LP1912'
# new row 96 (freshly scraped)
$ws.Cells.Item(96, 1).Value = '07:38:30'
$ws.Cells.Item(96, 2).Value = '09:29'
$ws.Cells.Item(96, 3).Value = '10_OLMOS'
$ws.Cells.Item(96, 4).Value = 111
$ws.Cells.Item(96, 5).Value = 'LP1912'
# new row 97 (freshly scraped)
$ws.Cells.Item(97, 1).Value = '07:38:30'
$ws.Cells.Item(97, 2).Value = '09:34'
$ws.Cells.Item(97, 3).Value = '15_ABASTO'
$ws.Cells.Item(97, 4).Value = 116
$ws.Cells.Item(97, 5).Value = 'LP1912'

# ----- LP1912-215 (6 changed rows) -----
$ws = $wb.Worksheets.Item('LP1912-215')
$ws.Cells.Item(2, 1).Value = 'Última actualización: 07:38:30'
$ws.Cells.Item(3, 1).Value = 'Total filas: 18'
$ws.Cells.Item(19, 1).Value = '07:38:30'
$ws.Cells.Item(19, 4).Value = 56
$ws.Cells.Item(21, 1).Value = '07:38:30'
$ws.Cells.Item(21, 2).Value = '08:48'
$ws.Cells.Item(21, 3).Value = '215A_EL PATO'
$ws.Cells.Item(21, 4).Value = 70
# new row 22 (freshly scraped)
$ws.Cells.Item(22, 1).Value = '07:38:30'
$ws.Cells.Item(22, 2).Value = '08:59'
$ws.Cells.Item(22, 3).Value = '215B_EL PATO'
$ws.Cells.Item(22, 4).Value = 81
$ws.Cells.Item(22, 5).Value = 'LP1912'
# new row 23 (freshly scraped)
$ws.Cells.Item(23, 1).Value = '07:38:30'
$ws.Cells.Item(23, 2).Value = '09:19'
$ws.Cells.Item(23, 3).Value = '215_EL PELIGRO'
$ws.Cells.Item(23, 4).Value = 101
$ws.Cells.Item(23, 5).Value = 'LP1912'

# ----- 6203-6173 (5 changed rows) -----
$ws = $wb.Worksheets.Item('6203-6173')
$ws.Cells.Item(2, 1).Value = 'Última actualización: 07:38:30'
$ws.Cells.Item(3, 1).Value = 'Total filas: 7'
$ws.Cells.Item(9, 1).Value = '07:38:30'
$ws.Cells.Item(9, 4).Value = 32
$ws.Cells.Item(11, 1).Value = '07:38:30'
$ws.Cells.Item(11, 4).Value = 45
# new row 12 (freshly scraped)
$ws.Cells.Item(12, 1).Value = '07:38:30'
$ws.Cells.Item(12, 2).Value = '08:52'
$ws.Cells.Item(12, 3).Value = '215A_LA PLATA'
$ws.Cells.Item(12, 4).Value = 74
$ws.Cells.Item(12, 5).Value = 'L6173'

